# Applies the "Fruta / hortaliza, semanal" update: the D (Fecha), J (Volumen),
# K (Precio mínimo), L (Precio máximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) columns for rows 2-23 are reshuffled among the rows
# (a permutation of the existing weekly records), while every other column
# stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the current (pre-edit) values for the columns that move.
$orig = @{}
for ($r = 2; $r -le 23; $r++) {
    $orig[$r] = @(
        $ws.Cells.Item($r, 4).Value2,   # D - Fecha
        $ws.Cells.Item($r, 10).Value2,  # J - Volumen
        $ws.Cells.Item($r, 11).Value2,  # K - Precio minimo
        $ws.Cells.Item($r, 12).Value2,  # L - Precio maximo
        $ws.Cells.Item($r, 13).Value2,  # M - Precio promedio ponderado
        $ws.Cells.Item($r, 16).Value2   # P - Precio $/Kg
    )
}

# Mapping: new row -> source row (data that row now contains came from this row before the edit).
$sourceRow = @{
    2 = 7; 3 = 21; 4 = 17; 5 = 20; 6 = 13; 7 = 3; 8 = 16; 9 = 15; 10 = 6;
    11 = 10; 12 = 18; 13 = 22; 14 = 19; 15 = 23; 16 = 9; 17 = 14; 18 = 5;
    19 = 8; 20 = 12; 21 = 11; 22 = 4; 23 = 2
}

foreach ($r in $sourceRow.Keys) {
    $src = $sourceRow[$r]
    $vals = $orig[$src]

    $ws.Cells.Item($r, 4).Value2 = $vals[0]
    $ws.Cells.Item($r, 10).Value2 = $vals[1]
    $ws.Cells.Item($r, 11).Value2 = $vals[2]
    $ws.Cells.Item($r, 12).Value2 = $vals[3]
    $ws.Cells.Item($r, 13).Value2 = $vals[4]
    $ws.Cells.Item($r, 16).Value2 = $vals[5]
}
